# Insert a new weekly data row at row 773 (pushing the existing rows 773-839
# down to 774-840), then populate the new row with the latest week's prices.
# Everything else (A/B/C/E/F/G/H/I/J/N/O/Q/R) for this new row mirrors what
# used to be in row 773, only the date (D) and the price columns (K/L/M/P)
# change to the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(773).Insert()

$ws.Cells.Item(773, 1).Value = 9
$ws.Cells.Item(773, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(773, 3).Value = "Metropolitana"
$ws.Cells.Item(773, 4).Value = 45166
$ws.Cells.Item(773, 5).Value = 13
$ws.Cells.Item(773, 6).Value = 100112031
$ws.Cells.Item(773, 7).Value = "Poroto verde"
$ws.Cells.Item(773, 8).Value = "Magnum"
$ws.Cells.Item(773, 9).Value = "Primera"
$ws.Cells.Item(773, 10).Value = 52
$ws.Cells.Item(773, 11).Value = 15000
$ws.Cells.Item(773, 12).Value = 17000
$ws.Cells.Item(773, 13).Value = 16000
$ws.Cells.Item(773, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(773, 15).Value = "Perú"
$ws.Cells.Item(773, 16).Value = 640
$ws.Cells.Item(773, 17).Value = 25
$ws.Cells.Item(773, 18).Value = "Hortaliza"
